$d = $word.ActiveDocument

function Replace-ParagraphXml($paragraph, [string]$innerXml) {
    $start = $paragraph.Range.Start
    $end   = $paragraph.Range.End
    $r = $d.Range($start, $end)
    $payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($payload) | Out-Null
}

# --- Paragraph 1 (title): "SnirtLabs Releases All New, Lightweight Document
# Collaboration System " -> bookmark + proofErr-wrapped "SnirtLabs" + new tagline
$p1Inner = @'
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="26B8FF"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="26B8FF"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>SnirtLabs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="26B8FF"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> Introduces a New, Ultra-lightweight Document Collaboration System </w:t></w:r>
'@
Replace-ParagraphXml $d.Paragraphs(1) $p1Inner

# --- Paragraph 4 (date/location + announcement lede). Keep the pPr and the
# unchanged tail (the <w:br/> runs + "When it comes..." text) intact, only
# rewrite the head portion.
$p4Inner = @'
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>SAN FRANCISCO, Calif. (August 3, 2014)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>SnirtLabs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> today announced the availability of a technical preview of Snirt its soon to be released lightweight document collaboration tool. </w:t></w:r><w:r w:rsidR="00B102AB" w:rsidRPr="00B102AB"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:br/></w:r><w:r w:rsidR="00B102AB" w:rsidRPr="00B102AB"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:br/><w:t>When it comes to documents we’ve all heard “Can you get me the latest version?” or “Which version is that?” But in the world of software development those questions were answered years ago through</w:t></w:r><w:r w:rsidR="00B102AB" w:rsidRPr="00B102AB"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> an invention called version control. It keeps track of even the smallest changes — over millions of files — without breaking a sweat. But where is version control for Microsoft Word documents, or my team’s PowerPoint project, or even a huge folder of Excel spreadsheets?   </w:t></w:r>
'@
Replace-ParagraphXml $d.Paragraphs(4) $p4Inner

# --- Paragraph 8 (the Nigel Chanter quote). Wrap the "SnirtLabs" mention in
# proofErr tags, split the run around it, keep the rest of the text identical.
$p8Inner = @'
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">"People don't want tools to get in the way of how they like to work and there are hundreds of millions of Microsoft Office users around the world attaching documents to emails as a basic means of collaboration, " said Nigel Chanter co-founder and CEO at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>SnirtLabs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>. "Our approach is to bring the power of collaboration directly inside the tools people already use, because it lets people work the way they want to."</w:t></w:r>
'@
Replace-ParagraphXml $d.Paragraphs(8) $p8Inner

# --- Paragraph 10 ("About SnirtLabs" heading) -> "About " + proofErr-wrapped "SnirtLabs"
$p10Inner = @'
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:b/></w:rPr><w:t xml:space="preserve">About </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:b/></w:rPr><w:t>SnirtLabs</w:t></w:r><w:proofErr w:type="spellEnd"/>
'@
Replace-ParagraphXml $d.Paragraphs(10) $p10Inner

# --- Paragraph 11 ("SnirtLabs is the developer..." body) -> wrap "SnirtLabs"
# (x3), "digial" and "Powerpoint" with proofErr spell-check markers.
$p11Inner = @'
<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>SnirtLabs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> is the developer of Snirt the new, ultra-lightweight document collaboration system. Snirt provides easy file sharing for all types of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>digial</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> content, effortless collaboration, secure storage and a visible version history of every document and file. A powerful add-in to Microsoft Word, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Powerpoint</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and Excel takes Snirt directly to all your work. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>SnirtLabs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> is headquartered in San Francisco, California. For more information visit www.Perforce.com.</w:t></w:r>
'@
Replace-ParagraphXml $d.Paragraphs(11) $p11Inner

Write-Host "Done applying edits"
